$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper fragments (OOXML wrapped in the WordOpenXML "pkg:package" shape
# that Range.InsertXML expects). Inserting at a COLLAPSED (zero-length)
# range keeps the new run(s) distinct from their neighbours instead of
# merging formatting-identical runs together.
# ---------------------------------------------------------------------
function New-RunsXml($inner) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $inner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$xmlExhibitOnly = New-RunsXml('<w:r><w:t>EXHIBIT</w:t></w:r>')
$xmlExhibitPlusParen = New-RunsXml('<w:r><w:t>EXHIBIT</w:t></w:r><w:r><w:t>(</w:t></w:r>')
$xmlExhibitBookmarkSpaceParen = New-RunsXml('<w:r><w:t>EXHIBIT</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(</w:t></w:r>')

# ---------------------------------------------------------------------
# 1. Drop the stray _GoBack bookmark that was sitting on the "SAMPLE 2"
#    heading paragraph.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. "UNF" sample-2 relation: prefix the bare "(" run with an "EXHIBIT"
#    run -> "EXHIBIT(ART_NO,ART_TITLE,ARTIST_CODE,...".
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(54)
$r1 = $p1.Range
$ins1 = $d.Range($r1.Start, $r1.Start)
$ins1.InsertXML($xmlExhibitOnly)

# ---------------------------------------------------------------------
# 3. "1NF" sample-2 relation: same prefixing of the leading "(" run.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(57)
$r2 = $p2.Range
$ins2 = $d.Range($r2.Start, $r2.Start)
$ins2.InsertXML($xmlExhibitOnly)

# ---------------------------------------------------------------------
# 4. 2NF EXHIBITION( ... ) line -> split the "EXHIBITION(" run into
#    separate "EXHIBIT" and "(" runs.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(67)
$r3 = $p3.Range
$del3 = $d.Range($r3.Start, $r3.Start + 11)
$del3.Text = ""
$ins3 = $d.Range($r3.Start, $r3.Start)
$ins3.InsertXML($xmlExhibitPlusParen)

# ---------------------------------------------------------------------
# 5. 3NF EXHIBITION( ... ) line -> "EXHIBIT" + the relocated _GoBack
#    bookmark + a literal space run + the "(" run.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(77)
$r4 = $p4.Range
$del4 = $d.Range($r4.Start, $r4.Start + 11)
$del4.Text = ""
$ins4 = $d.Range($r4.Start, $r4.Start)
$ins4.InsertXML($xmlExhibitBookmarkSpaceParen)

# ---------------------------------------------------------------------
# 6. Attribute-synthesis EXHIBITION(...) line -> rename the standalone
#    "EXHIBITION" run to "EXHIBIT" (the "(" is already its own run).
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(98)
$r5 = $p5.Range
$del5 = $d.Range($r5.Start, $r5.Start + 10)
$del5.Text = ""
$ins5 = $d.Range($r5.Start, $r5.Start)
$ins5.InsertXML($xmlExhibitOnly)

# ---------------------------------------------------------------------
# 7. Register the (until now latent) "Balloon Text" paragraph style and
#    its linked "Balloon Text Char" character style, matching what Word
#    stamps into styles.xml the first time a comment/tracked-change
#    balloon is touched during the session.
# ---------------------------------------------------------------------
$balloon = $d.Styles.Add("Balloon Text", 1)
$balloon.BaseStyle = "Normal"
$balloon.LinkStyle = "BalloonTextChar"
$balloon.Priority = 99
$balloon.UnhideWhenUsed = $true
$balloon.Font.Name = "Times New Roman"
$balloon.Font.NameAscii = "Times New Roman"
$balloon.Font.NameBi = "Times New Roman"
$balloon.Font.Size = 9
$balloon.Font.SizeBi = 9

$balloonChar = $d.Styles.Add("Balloon Text Char", 2)
$balloonChar.BaseStyle = "DefaultParagraphFont"
$balloonChar.LinkStyle = "BalloonText"
$balloonChar.Priority = 99
$balloonChar.Font.Name = "Times New Roman"
$balloonChar.Font.NameAscii = "Times New Roman"
$balloonChar.Font.NameBi = "Times New Roman"
$balloonChar.Font.Size = 9
$balloonChar.Font.SizeBi = 9

Write-Host "Done"
